$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update department (C2): was the long faculty name, now shortened to "Community Services"
$ws.Range("C2").Value = "Community Services"

# Clear promotionValidity (R2): remove the expired promotion text, leaving the cell blank
$ws.Range("R2").Value = ""
